# Update "Metadata" sheet: refresh the "Last Updated" timestamp
$wbMeta = $excel.ActiveWorkbook
$wsMeta = $wbMeta.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 12:40 PM"

# Update "Stock List" sheet: new snapshot of stock data.
# A new row (CAPTRU-RE1) is inserted at the top of the data (row 2),
# pushing all following rows down by one, and the last row (row 76,
# previously TRAVELFOOD) drops off the bottom of the list.
$ws = $wbMeta.Worksheets.Item("Stock List")

$lastDataRow = 75   # rows 2..75 shift down into rows 3..76
for ($r = $lastDataRow; $r -ge 2; $r--) {
    $nr = $r + 1
    $ws.Cells.Item($nr, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($nr, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($nr, 4).Value = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($nr, 5).Value = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($nr, 8).Value = $ws.Cells.Item($r, 8).Value2
}

# New row 2 values (A, F, G columns are unchanged: icon / N/A / N/A)
$ws.Cells.Item(2, 2).Value = "CAPTRU-RE1"
$ws.Cells.Item(2, 3).Value = "CAPTRU-RE1"
$ws.Cells.Item(2, 4).Value = 5.67
$ws.Cells.Item(2, 5).Value = -11.9565
